# Auto-generated edit script: refresh cached market-price columns (H-N)
# on each job sheet, per the scheduled-runner data pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 238.26666
$ws.Range("I33").Value = 131.63637
$ws.Range("K33").Value = 131.63637
$ws.Range("M33").Value = 97.36363
$ws.Range("H98").Value = 4925.778
$ws.Range("I98").Value = 5457.2666
$ws.Range("J98").Value = 2268.3333
$ws.Range("K98").Value = 5457.2666
$ws.Range("L98").Value = 2268.3333
$ws.Range("M98").Value = -3959.2666
$ws.Range("N98").Value = -5264.3333
$ws.Range("H106").Value = 11017
$ws.Range("I106").Value = 11402.462
$ws.Range("K106").Value = 11402.462
$ws.Range("M106").Value = -10771.462
$ws.Range("H111").Value = 2857.0715
$ws.Range("I111").Value = 2334.3333
$ws.Range("K111").Value = 7002.999899999999
$ws.Range("M111").Value = -3935.999899999999
$ws.Range("H113").Value = 2758.1333
$ws.Range("I113").Value = 2547.7
$ws.Range("K113").Value = 2547.7
$ws.Range("M113").Value = 706.3000000000002
$ws.Range("H122").Value = 4925.778
$ws.Range("I122").Value = 5457.2666
$ws.Range("J122").Value = 2268.3333
$ws.Range("K122").Value = 16371.7998
$ws.Range("L122").Value = 6804.999899999999
$ws.Range("M122").Value = -13921.7998
$ws.Range("N122").Value = -11704.9999
$ws.Range("H132").Value = 10426585
$ws.Range("I132").Value = 20842694
$ws.Range("J132").Value = 10476.0625
$ws.Range("K132").Value = 62528082
$ws.Range("L132").Value = 31428.1875
$ws.Range("M132").Value = -62525552
$ws.Range("N132").Value = -36488.1875
$ws.Range("H135").Value = 560.1081
$ws.Range("I135").Value = 227.17647
$ws.Range("K135").Value = 2044.58823
$ws.Range("M135").Value = 490.4117700000002
$ws.Range("H137").Value = 1099.4667
$ws.Range("I137").Value = 910.925
$ws.Range("J137").Value = 1314.9429
$ws.Range("K137").Value = 2732.775
$ws.Range("L137").Value = 3944.8287
$ws.Range("M137").Value = -182.7749999999996
$ws.Range("N137").Value = -9044.8287
$ws.Range("H139").Value = 34840
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 34840
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 34840
$ws.Range("M139").ClearContents()
$ws.Range("N139").Value = -45120
$ws.Range("H141").Value = 530.5217
$ws.Range("I141").Value = 530.5217
$ws.Range("K141").Value = 1591.5651
$ws.Range("M141").Value = 3588.4349

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 15633.714
$ws.Range("I2").Value = 1408.25
$ws.Range("J2").Value = 34601
$ws.Range("K2").Value = 1408.25
$ws.Range("L2").Value = 34601
$ws.Range("M2").Value = -1295.25
$ws.Range("N2").Value = -34827
$ws.Range("H61").Value = 22223150
$ws.Range("I61").Value = 24391070
$ws.Range("J61").Value = 1953.5
$ws.Range("K61").Value = 24391070
$ws.Range("L61").Value = 1953.5
$ws.Range("M61").Value = -24390858
$ws.Range("N61").Value = -2377.5
$ws.Range("H74").Value = 815.72974
$ws.Range("I74").Value = 710.82355
$ws.Range("J74").Value = 2004.6666
$ws.Range("K74").Value = 710.82355
$ws.Range("L74").Value = 2004.6666
$ws.Range("M74").Value = 163.17645
$ws.Range("N74").Value = -3752.6666
$ws.Range("H77").Value = 815.72974
$ws.Range("I77").Value = 710.82355
$ws.Range("J77").Value = 2004.6666
$ws.Range("K77").Value = 3554.11775
$ws.Range("L77").Value = 10023.333
$ws.Range("M77").Value = 813.8822500000001
$ws.Range("N77").Value = -18759.333
$ws.Range("H116").Value = 15633.714
$ws.Range("I116").Value = 1408.25
$ws.Range("J116").Value = 34601
$ws.Range("K116").Value = 1408.25
$ws.Range("L116").Value = 34601
$ws.Range("M116").Value = 885.75
$ws.Range("N116").Value = -39189
$ws.Range("H122").Value = 2823.4443
$ws.Range("I122").Value = 2738.875
$ws.Range("K122").Value = 8216.625
$ws.Range("M122").Value = -5766.625
$ws.Range("H132").Value = 2451.1
$ws.Range("I132").Value = 2371.7932
$ws.Range("J132").Value = 2660.182
$ws.Range("K132").Value = 7115.3796
$ws.Range("L132").Value = 7980.545999999999
$ws.Range("M132").Value = -4585.3796
$ws.Range("N132").Value = -13040.546
$ws.Range("H136").Value = 22223150
$ws.Range("I136").Value = 24391070
$ws.Range("J136").Value = 1953.5
$ws.Range("K136").Value = 73173210
$ws.Range("L136").Value = 5860.5
$ws.Range("M136").Value = -73170660
$ws.Range("N136").Value = -10960.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 15633.714
$ws.Range("I3").Value = 1408.25
$ws.Range("J3").Value = 34601
$ws.Range("K3").Value = 1408.25
$ws.Range("L3").Value = 34601
$ws.Range("M3").Value = -1294.25
$ws.Range("N3").Value = -34829
$ws.Range("H107").Value = 2081
$ws.Range("I107").Value = 1571.1428
$ws.Range("J107").Value = 3270.6667
$ws.Range("K107").Value = 1571.1428
$ws.Range("L107").Value = 3270.6667
$ws.Range("M107").Value = 348.8571999999999
$ws.Range("N107").Value = -7110.6667
$ws.Range("H134").Value = 4581.8223
$ws.Range("I134").Value = 1784.3529
$ws.Range("J134").Value = 13228.546
$ws.Range("K134").Value = 5353.0587
$ws.Range("L134").Value = 39685.638
$ws.Range("M134").Value = -2818.0587
$ws.Range("N134").Value = -44755.638

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1916.4642
$ws.Range("I31").Value = 2058.5652
$ws.Range("J31").Value = 1262.8
$ws.Range("K31").Value = 2058.5652
$ws.Range("L31").Value = 1262.8
$ws.Range("M31").Value = -1763.5652
$ws.Range("N31").Value = -1852.8
$ws.Range("H34").Value = 1916.4642
$ws.Range("I34").Value = 2058.5652
$ws.Range("J34").Value = 1262.8
$ws.Range("K34").Value = 2058.5652
$ws.Range("L34").Value = 1262.8
$ws.Range("M34").Value = -1856.5652
$ws.Range("N34").Value = -1666.8
$ws.Range("H58").Value = 810.0857
$ws.Range("I58").Value = 719.96295
$ws.Range("K58").Value = 719.96295
$ws.Range("M58").Value = -516.96295
$ws.Range("H105").Value = 676
$ws.Range("I105").Value = 607.3333
$ws.Range("K105").Value = 607.3333
$ws.Range("M105").Value = 1139.6667
$ws.Range("H132").Value = 2719.9062
$ws.Range("I132").Value = 2434.652
$ws.Range("K132").Value = 7303.956
$ws.Range("M132").Value = -4773.956
$ws.Range("H136").Value = 810.0857
$ws.Range("I136").Value = 719.96295
$ws.Range("K136").Value = 2159.88885
$ws.Range("M136").Value = 390.1111500000002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 1388.4
$ws.Range("I32").Value = 797.3333
$ws.Range("J32").Value = 2275
$ws.Range("K32").Value = 2391.9999
$ws.Range("L32").Value = 6825
$ws.Range("M32").Value = -2108.9999
$ws.Range("N32").Value = -7391
$ws.Range("H134").Value = 3056.7932
$ws.Range("I134").Value = 1307.1875
$ws.Range("J134").Value = 5210.154
$ws.Range("K134").Value = 3921.5625
$ws.Range("L134").Value = 15630.462
$ws.Range("M134").Value = 1148.4375
$ws.Range("N134").Value = -25770.462
$ws.Range("H136").Value = 2490.6875
$ws.Range("J136").Value = 3370.111
$ws.Range("L136").Value = 10110.333
$ws.Range("N136").Value = -20310.333

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2350.6
$ws.Range("I122").Value = 1688.375
$ws.Range("K122").Value = 5065.125
$ws.Range("M122").Value = -2615.125
$ws.Range("H132").Value = 1930.2759
$ws.Range("I132").Value = 1518.35
$ws.Range("K132").Value = 4555.049999999999
$ws.Range("M132").Value = -2025.049999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5650.5
$ws.Range("I40").Value = 2725.75
$ws.Range("K40").Value = 2725.75
$ws.Range("M40").Value = -2589.75
$ws.Range("H122").Value = 14709583
$ws.Range("I122").Value = 20834568
$ws.Range("J122").Value = 9619.799999999999
$ws.Range("K122").Value = 62503704
$ws.Range("L122").Value = 28859.4
$ws.Range("M122").Value = -62501254
$ws.Range("N122").Value = -33759.39999999999
$ws.Range("H132").Value = 20676.094
$ws.Range("I132").Value = 1248.9656
$ws.Range("J132").Value = 44150.543
$ws.Range("K132").Value = 3746.8968
$ws.Range("L132").Value = 132451.629
$ws.Range("M132").Value = -1216.8968
$ws.Range("N132").Value = -137511.629
$ws.Range("H136").Value = 1165.2258
$ws.Range("I136").Value = 1078.5927
$ws.Range("J136").Value = 1750
$ws.Range("K136").Value = 3235.7781
$ws.Range("L136").Value = 5250
$ws.Range("M136").Value = -685.7780999999995
$ws.Range("N136").Value = -10350

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 56821920
$ws.Range("I122").Value = 59527484
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 178582452
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -178580002
$ws.Range("N122").Value = -19900
$ws.Range("H132").Value = 2424.8408
$ws.Range("I132").Value = 2508.1943
$ws.Range("K132").Value = 7524.5829
$ws.Range("M132").Value = -4994.5829
$ws.Range("H136").Value = 549.9388
$ws.Range("I136").Value = 398.13513
$ws.Range("J136").Value = 1018
$ws.Range("K136").Value = 1194.40539
$ws.Range("L136").Value = 3054
$ws.Range("M136").Value = 1355.59461
$ws.Range("N136").Value = -8154
